$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-27 Wednesday", "2025-08-28 Thursday"),
    @("148÷5=29, 3", "809÷9=89, 8"),
    @("201÷5=40, 1", "271÷7=38, 5"),
    @("429÷6=71, 3", "841÷8=105, 1"),
    @("918÷8=114, 6", "303÷6=50, 3"),
    @("338÷8=42, 2", "886÷4=221, 2"),
    @("327÷7=46, 5", "845÷9=93, 8"),
    @("526÷8=65, 6", "324÷3=108, 0"),
    @("916÷9=101, 7", "432÷5=86, 2"),
    @("207÷2=103, 1", "880÷2=440, 0"),
    @("203÷9=22, 5", "622÷6=103, 4"),
    @("592÷3=197, 1", "482÷9=53, 5"),
    @("724÷8=90, 4", "984÷7=140, 4"),
    @("141÷7=20, 1", "543÷7=77, 4"),
    @("654÷3=218, 0", "176÷8=22, 0"),
    @("759÷6=126, 3", "908÷7=129, 5"),
    @("602÷2=301, 0", "546÷4=136, 2"),
    @("762÷9=84, 6", "517÷8=64, 5"),
    @("651÷3=217, 0", "782÷4=195, 2"),
    @("587÷8=73, 3", "527÷8=65, 7"),
    @("233÷6=38, 5", "776÷6=129, 2"),
    @("593÷5=118, 3", "289÷4=72, 1"),
    @("478÷6=79, 4", "567÷9=63, 0"),
    @("697÷3=232, 1", "204÷4=51, 0"),
    @("461÷2=230, 1", "338÷5=67, 3"),
    @("945÷6=157, 3", "494÷7=70, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
